# Update "想去人数" (number of interested attendees) figures in both the
# "展览" and "全部类型" sheets, which hold duplicate data tables.
$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 129
    $ws.Range("F5").Value = 25
    $ws.Range("F6").Value = 442
    $ws.Range("F9").Value = 567
    $ws.Range("F10").Value = 404
}
